$d = $word.ActiveDocument

# Word constant used by Find.Execute (Wrap:=wdFindContinue)
$wdFindContinue = 1
$wdReplaceAll = 2

# ---------------------------------------------------------------------------
# 1. Point the "(" + email ")" hyperlink at the new mailbox and update the
#    visible email text to match. Do this before touching the neighboring
#    name text so this paragraph only needs one more edit afterwards (the
#    run-splitting trick below relies on being the *last* edit made to the
#    paragraph it lives in).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("salz2@illinois.edu", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "amshim2@illinois.edu", $wdReplaceAll)

foreach ($hl in $d.Hyperlinks) {
    if ($hl.Address -like "mailto:salz2@illinois.edu*") {
        $hl.Address = "mailto:amshim2@illinois.edu"
    }
}

# ---------------------------------------------------------------------------
# 2. Move the "_GoBack" last-edit bookmark from the end of the document to
#    sit inside the word "members" (after "board m", before "embers") in the
#    "contact current board members" bullet, matching where the edit
#    session's cursor last sat. Adding a bookmark with an existing name
#    relocates it (and drops the old one), just like real Word.
# ---------------------------------------------------------------------------
$membersRng = $d.Content.Duplicate
$membersRng.Find.Execute("board members", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
$splitPoint = $membersRng.Start + 7
$goBackRng = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $goBackRng)

# ---------------------------------------------------------------------------
# 3. Swap the contact name "Brady Salz" -> "Anselmo Shim".
#    Only the name itself is replaced (not the trailing space) so the
#    untouched trailing-space text keeps living in its own run, the same
#    way Word splits a run when only part of it is retyped. This is the
#    last edit in the "Submit the electronic application..." paragraph so
#    the split actually survives the save.
# ---------------------------------------------------------------------------
$nameRng = $d.Content.Duplicate
$nameRng.Find.Execute("Brady Salz ", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
$nameStart = $nameRng.Start
$nameEnd = $nameRng.End

$justName = $d.Range($nameStart, $nameEnd - 1)
$justName.Text = "Anselmo Shim"

# Force the remaining trailing-space text to stay in its own run instead of
# being re-merged with the freshly typed "Anselmo Shim" run.
$spaceRng = $d.Range($nameStart + 12, $nameStart + 13)
$spaceRng.Bold = 1
$spaceRng.Bold = 0
